# Fruta / hortaliza, semanal
# Appends three new weekly price rows (82-84) for
# "Comercializadora del Agro de Limarí" - Durazno, variedad Andross.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 82; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44595; E = 4; F = "Fruta"; G = 100103; H = "Frutos de hueso (carozo)"; I = 100103004; J = "Durazno"; K = "Andross"; L = "Especial"; M = 200;  N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 1417; T = 18 },
    @{ Row = 83; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44595; E = 4; F = "Fruta"; G = 100103; H = "Frutos de hueso (carozo)"; I = 100103004; J = "Durazno"; K = "Andross"; L = "Primera";  M = 240;  N = 20000; O = 21000; P = 20500; Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 1139; T = 18 },
    @{ Row = 84; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44595; E = 4; F = "Fruta"; G = 100103; H = "Frutos de hueso (carozo)"; I = 100103004; J = "Durazno"; K = "Andross"; L = "Segunda";  M = 300;  N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos empedrada"; R = "Región de O'Higgins"; S = 861;  T = 18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
